{"js": "// Replace the three-digit-division problems in the table with their\n// updated values. Each original expression (e.g. \"378\u00f73=\") is\n// unique in the document, so body.search(...) + insertText(\"Replace\")\n// swaps the text of the matching run while preserving its formatting.\nconst replacements = [\n  [\"378\u00f73=\", \"482\u00f79=\"],\n  [\"791\u00f74=\", \"322\u00f74=\"],\n  [\"731\u00f76=\", \"449\u00f73=\"],\n  [\"899\u00f75=\", \"218\u00f74=\"],\n  [\"348\u00f73=\", \"873\u00f75=\"],\n  [\"584\u00f74=\", \"884\u00f74=\"],\n  [\"617\u00f79=\", \"732\u00f76=\"],\n  [\"960\u00f73=\", \"827\u00f77=\"],\n  [\"389\u00f77=\", \"242\u00f77=\"],\n  [\"904\u00f76=\", \"476\u00f77=\"],\n  [\"469\u00f75=\", \"666\u00f77=\"],\n  [\"751\u00f76=\", \"208\u00f73=\"],\n  [\"853\u00f75=\", \"587\u00f72=\"],\n  [\"296\u00f77=\", \"241\u00f78=\"],\n  [\"133\u00f72=\", \"739\u00f77=\"],\n  [\"626\u00f75=\", \"215\u00f73=\"],\n  [\"469\u00f72=\", \"226\u00f75=\"],\n  [\"204\u00f76=\", \"633\u00f73=\"],\n  [\"586\u00f78=\", \"403\u00f73=\"],\n  [\"974\u00f72=\", \"388\u00f74=\"],\n  [\"757\u00f74=\", \"641\u00f74=\"],\n  [\"117\u00f79=\", \"621\u00f75=\"],\n  [\"607\u00f78=\", \"930\u00f75=\"],\n  [\"360\u00f75=\", \"590\u00f79=\"],\n  [\"923\u00f78=\", \"615\u00f73=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit-division problems in the table to their new\n# values. Each original expression (e.g. \"378\u00f73=\") is unique in\n# the document, so Find/Replace with wdReplaceAll swaps exactly one run\n# of text per pair while leaving all other formatting untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"378\u00f73=\", \"482\u00f79=\"),\n    @(\"791\u00f74=\", \"322\u00f74=\"),\n    @(\"731\u00f76=\", \"449\u00f73=\"),\n    @(\"899\u00f75=\", \"218\u00f74=\"),\n    @(\"348\u00f73=\", \"873\u00f75=\"),\n    @(\"584\u00f74=\", \"884\u00f74=\"),\n    @(\"617\u00f79=\", \"732\u00f76=\"),\n    @(\"960\u00f73=\", \"827\u00f77=\"),\n    @(\"389\u00f77=\", \"242\u00f77=\"),\n    @(\"904\u00f76=\", \"476\u00f77=\"),\n    @(\"469\u00f75=\", \"666\u00f77=\"),\n    @(\"751\u00f76=\", \"208\u00f73=\"),\n    @(\"853\u00f75=\", \"587\u00f72=\"),\n    @(\"296\u00f77=\", \"241\u00f78=\"),\n    @(\"133\u00f72=\", \"739\u00f77=\"),\n    @(\"626\u00f75=\", \"215\u00f73=\"),\n    @(\"469\u00f72=\", \"226\u00f75=\"),\n    @(\"204\u00f76=\", \"633\u00f73=\"),\n    @(\"586\u00f78=\", \"403\u00f73=\"),\n    @(\"974\u00f72=\", \"388\u00f74=\"),\n    @(\"757\u00f74=\", \"641\u00f74=\"),\n    @(\"117\u00f79=\", \"621\u00f75=\"),\n    @(\"607\u00f78=\", \"930\u00f75=\"),\n    @(\"360\u00f75=\", \"590\u00f79=\"),\n    @(\"923\u00f78=\", \"615\u00f73=\"),\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
